# Update time zone labels from UTC to UTC + 2 in the data description sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_description")

$ws.Range("C2").Value  = "Date and time of sampling (UTC + 2)"
$ws.Range("C32").Value = "Start time of light measurement  (UTC + 2)"
$ws.Range("C33").Value = "End time of light measurement (UTC + 2)"
$ws.Range("C47").Value = "Time in seconds for dark measurement  (UTC + 2)"
$ws.Range("C48").Value = "PAR value of dark measurement  (UTC + 2)"
$ws.Range("C57").Value = "Time of sampling  (UTC + 2)"

# Move the active selection from C33 to C4, matching the saved view state.
$ws.Activate()
$ws.Range("C4").Select()
